$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Factors")
$lo = $ws.ListObjects.Item("Table13")

# Insert a new column before column C (Location is B, Summary is currently C) -
# this shifts all cell data (and the existing table columns' underlying cells)
# one column to the right.
$ws.Range("C:C").Insert()

# Grow the table definition so it once again spans the (now wider) data range.
$lo.Resize($ws.Range("A2:G1048576"))

# The ListObject's cached column names don't automatically resync with the
# shifted header cells, so re-assert them left-to-right to match the actual
# (shifted) header text, finishing with the brand-new header in column C.
$ws.Cells.Item(2, 7).Value = $ws.Cells.Item(2, 7).Value()
$ws.Cells.Item(2, 6).Value = $ws.Cells.Item(2, 6).Value()
$ws.Cells.Item(2, 5).Value = $ws.Cells.Item(2, 5).Value()
$ws.Cells.Item(2, 4).Value = $ws.Cells.Item(2, 4).Value()
$ws.Range("C2").Value = "Objective"

$ws.Columns("C").ColumnWidth = 14

# Keep the same cell selection that was active before editing.
$ws.Range("C3").Select()

$gaps = $wb.Worksheets.Item("Gaps")
$gaps.Range("C2").Select()
